# Daily update at 8 AM UTC
# Appends the next day's row (row 69) to the "Wins Over Time" tracker sheet
# and moves the "last row" date formatting down from row 68 to row 69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 68 is no longer the final row, so its date cell reverts to the
# regular date number format used by every other data row (same as A67).
$ws.Range("A68").NumberFormat = $ws.Range("A67").NumberFormat

# Append the new day's data as row 69.
$ws.Range("A69").Value = 45656
$ws.Range("B69").Value = 163
$ws.Range("C69").Value = 155
$ws.Range("D69").Value = 160

# The new last row takes on the distinct "final row" date number format
# that used to live on A68.
$ws.Range("A69").NumberFormat = "YYYY-MM-DD"
